$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Item Name"
$ws.Range("C1").Value = "Buy Price "
$ws.Range("D1").Value = "Sell Price"
$ws.Range("E1").Value = "Quantity"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "shark"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 11
$ws.Range("E2").Value = 12

$null = $ws.Range("F1").Select()
